$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellValues = @{
    "G2" = 10.12586833333333
    "H2" = 30.377605
    "I2" = 0.9311967029481902
    "J2" = 0.9311967029481902
    "M2" = 7.413580666666667
    "N2" = 22.240742
    "O2" = 0.05108888817597561
    "P2" = 0.05108888817597561
    "Q2" = 75.06894170921224
    "R2" = 675.6204753829101
    "S2" = 0.04757380422675727
    "T2" = 0.04757380422675727
    "G3" = 10.12586833333333
    "H3" = 30.377605
    "I3" = 0.9311967029481902
    "J3" = 0.9311967029481902
    "O3" = 0.0112127179963522
    "P3" = 0.0112127179963522
    "Q3" = 16.47573285937778
    "R3" = 148.2815957344
    "S3" = 0.01044124602929101
    "T3" = 0.01044124602929101
    "G4" = 10.12586833333333
    "H4" = 30.377605
    "I4" = 0.9311967029481902
    "J4" = 0.9311967029481902
    "O4" = 0.9376983938276722
    "P4" = 0.9376983938276722
    "Q4" = 1377.834370256918
    "R4" = 12400.50933231226
    "S4" = 0.873181652692142
    "T4" = 0.873181652692142
    "I5" = 0.009287810103293732
    "J5" = 0.009287810103293733
    "M5" = 7.413580666666667
    "N5" = 22.240742
    "O5" = 0.05108888817597561
    "P5" = 0.05108888817597561
    "Q5" = 0.7487419930106667
    "R5" = 6.738677937096
    "S5" = 0.00047450389176687
    "T5" = 0.00047450389176687
    "I6" = 0.009287810103293732
    "J6" = 0.009287810103293733
    "O6" = 0.0112127179963522
    "P6" = 0.0112127179963522
    "Q6" = 0.1643299182933333
    "S6" = 0.0001041415954919034
    "T6" = 0.0001041415954919034
    "I7" = 0.009287810103293732
    "J7" = 0.009287810103293733
    "O7" = 0.9376983938276722
    "P7" = 0.9376983938276722
    "Q7" = 13.74260018771733
    "S7" = 0.008709164616034958
    "T7" = 0.00870916461603496
    "I8" = 0.05951548694851595
    "J8" = 0.05951548694851596
    "M8" = 7.413580666666667
    "N8" = 22.240742
    "O8" = 0.05108888817597561
    "P8" = 0.05108888817597561
    "Q8" = 4.797874183175778
    "R8" = 43.180867648582
    "S8" = 0.003040580057451468
    "T8" = 0.003040580057451468
    "I9" = 0.05951548694851595
    "J9" = 0.05951548694851596
    "O9" = 0.0112127179963522
    "P9" = 0.0112127179963522
    "S9" = 0.0006673303715692894
    "T9" = 0.0006673303715692894
    "I10" = 0.05951548694851595
    "J10" = 0.05951548694851596
    "O10" = 0.9376983938276722
    "P10" = 0.9376983938276722
    "R10" = 792.552582054452
    "S10" = 0.05580757651949519
    "T10" = 0.0558075765194952
}

foreach ($cell in $cellValues.Keys) {
    $ws.Range($cell).Value = $cellValues[$cell]
}
